$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update Marking row (B11): number of right-answer marks
$ws.Range("B11").Value = 5

# Update Total row (B12): total marks obtained
$ws.Range("B12").Value = 115

# Update Total/Max display (E12): "corr/total" marks text
$ws.Range("E12").Value = "115/140"
